$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.234.25"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.989.50"
$ws.Range("E3").Value = "  +5.94%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.24"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5097"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4136"
$ws.Range("E8").Value = "  +4.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08652"
$ws.Range("E9").Value = "  +5.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.68"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.19"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "1.994.75"
$ws.Range("E13").Value = "  +6.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.484"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.374"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.87"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001111"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06546"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.68"
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.067"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("D23").Value = "30.299.96"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.55"
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.204"
$ws.Range("D26").Value = "2.225.15"
$ws.Range("E26").Value = "  +6.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.49"
$ws.Range("E27").Value = "  +6.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.06"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.354"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.43"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.128"
$ws.Range("E31").Value = "  +5.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1051"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.049"
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.813"
$ws.Range("E34").Value = "  +3.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.311"
$ws.Range("E35").Value = "  +11.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02476"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.379"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06509"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2191"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.873"
$ws.Range("E40").Value = "  +4.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6562"
$ws.Range("E41").Value = "  +4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.77"
$ws.Range("E42").Value = "  +4.44%  "
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.62"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6107"
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.195"
$ws.Range("E46").Value = "  +4.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.667"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.15"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.23"
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06872"
$ws.Range("E51").Value = "  +1.97%  "
